$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1555.619
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 1583.4
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 4750.200000000001
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -5086.200000000001
$ws.Range("H70").Value = 3907.318
$ws.Range("J70").Value = 2395.4285
$ws.Range("L70").Value = 7186.2855
$ws.Range("N70").Value = -7726.2855
$ws.Range("H73").Value = 3907.318
$ws.Range("J73").Value = 2395.4285
$ws.Range("L73").Value = 7186.2855
$ws.Range("N73").Value = -9058.2855
$ws.Range("H80").Value = 4661
$ws.Range("I80").Value = 6655.4287
$ws.Range("J80").Value = 2334.1667
$ws.Range("K80").Value = 19966.2861
$ws.Range("L80").Value = 7002.500100000001
$ws.Range("M80").Value = -18968.2861
$ws.Range("N80").Value = -8998.500100000001
$ws.Range("H83").Value = 4661
$ws.Range("I83").Value = 6655.4287
$ws.Range("J83").Value = 2334.1667
$ws.Range("K83").Value = 59898.85830000001
$ws.Range("L83").Value = 21007.5003
$ws.Range("M83").Value = -54906.85830000001
$ws.Range("N83").Value = -30991.5003
$ws.Range("H98").Value = 671.2857
$ws.Range("I98").Value = 671.2857
$ws.Range("K98").Value = 671.2857
$ws.Range("M98").Value = 826.7143
$ws.Range("H122").Value = 671.2857
$ws.Range("I122").Value = 671.2857
$ws.Range("K122").Value = 2013.8571
$ws.Range("M122").Value = 436.1428999999998
$ws.Range("H127").Value = 4381.5713
$ws.Range("I127").Value = 4381.5713
$ws.Range("K127").Value = 13144.7139
$ws.Range("M127").Value = -8184.713899999999
$ws.Range("H129").Value = 1029.7391
$ws.Range("I129").Value = 668
$ws.Range("K129").Value = 2004
$ws.Range("M129").Value = 2996
$ws.Range("H137").Value = 9797.333000000001
$ws.Range("I137").Value = 9797.333000000001
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 29391.999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -26841.999
$ws.Range("N137").Value = ""
$ws.Range("H138").Value = 4622.4
$ws.Range("I138").Value = 5584.7144
$ws.Range("K138").Value = 16754.1432
$ws.Range("M138").Value = -11614.1432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 225403.69
$ws.Range("I32").Value = 2212
$ws.Range("K32").Value = 2212
$ws.Range("M32").Value = -1925
$ws.Range("H63").Value = 6429.294
$ws.Range("I63").Value = 1824.5
$ws.Range("K63").Value = 1824.5
$ws.Range("M63").Value = -1138.5
$ws.Range("H66").Value = 6429.294
$ws.Range("I66").Value = 1824.5
$ws.Range("K66").Value = 9122.5
$ws.Range("M66").Value = -5690.5
$ws.Range("H97").Value = 507.1111
$ws.Range("I97").Value = 507.1111
$ws.Range("K97").Value = 507.1111
$ws.Range("M97").Value = -11.11110000000002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2509.4
$ws.Range("I86").Value = 2136.75
$ws.Range("K86").Value = 2136.75
$ws.Range("M86").Value = -1013.75
$ws.Range("H89").Value = 2509.4
$ws.Range("I89").Value = 2136.75
$ws.Range("K89").Value = 10683.75
$ws.Range("M89").Value = -5067.75
$ws.Range("H99").Value = 1484.85
$ws.Range("I99").Value = 1183.1666
$ws.Range("J99").Value = 4200
$ws.Range("K99").Value = 1183.1666
$ws.Range("L99").Value = 4200
$ws.Range("M99").Value = 314.8334
$ws.Range("N99").Value = -7196

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 4750
$ws.Range("I44").Value = 4750
$ws.Range("K44").Value = 4750
$ws.Range("M44").Value = -4308
$ws.Range("H62").Value = 85398.60000000001
$ws.Range("I62").Value = 6624.5
$ws.Range("K62").Value = 6624.5
$ws.Range("M62").Value = -6000.5
$ws.Range("H65").Value = 85398.60000000001
$ws.Range("I65").Value = 6624.5
$ws.Range("K65").Value = 33122.5
$ws.Range("M65").Value = -30002.5
$ws.Range("H99").Value = 11356.583
$ws.Range("I99").Value = 6592.85
$ws.Range("K99").Value = 6592.85
$ws.Range("M99").Value = -5094.85
$ws.Range("H126").Value = 11356.583
$ws.Range("I126").Value = 6592.85
$ws.Range("K126").Value = 19778.55
$ws.Range("M126").Value = -17308.55

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2421.037
$ws.Range("J131").Value = 2477.3726
$ws.Range("L131").Value = 7432.1178
$ws.Range("N131").Value = -17512.1178

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 31000
$ws.Range("J15").Value = 45000
$ws.Range("L15").Value = 45000
$ws.Range("N15").Value = -45576
$ws.Range("H81").Value = 31000
$ws.Range("J81").Value = 45000
$ws.Range("L81").Value = 45000
$ws.Range("N81").Value = -46996
$ws.Range("H84").Value = 31000
$ws.Range("J84").Value = 45000
$ws.Range("L84").Value = 135000
$ws.Range("N84").Value = -144984
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H122").Value = 48259.137
$ws.Range("I122").Value = 2147.5264
$ws.Range("K122").Value = 6442.5792
$ws.Range("M122").Value = -3992.5792

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5336.6665
$ws.Range("I7").Value = 5336.6665
$ws.Range("K7").Value = 5336.6665
$ws.Range("M7").Value = -5224.6665
$ws.Range("H55").Value = 746.86365
$ws.Range("I55").Value = 671.06665
$ws.Range("J55").Value = 909.2857
$ws.Range("K55").Value = 671.06665
$ws.Range("L55").Value = 909.2857
$ws.Range("M55").Value = -498.06665
$ws.Range("N55").Value = -1255.2857
$ws.Range("H122").Value = 2061.25
$ws.Range("I122").Value = 2061.25
$ws.Range("K122").Value = 6183.75
$ws.Range("M122").Value = -3733.75
$ws.Range("H126").Value = 5336.6665
$ws.Range("I126").Value = 5336.6665
$ws.Range("K126").Value = 16009.9995
$ws.Range("M126").Value = -13539.9995
$ws.Range("H134").Value = 44251.6
$ws.Range("J134").Value = 44251.6
$ws.Range("L134").Value = 44251.6
$ws.Range("N134").Value = -54391.6
$ws.Range("H136").Value = 8944
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 8944
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 26832
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -31932

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 49999.5
$ws.Range("I58").Value = 49999.5
$ws.Range("K58").Value = 49999.5
$ws.Range("M58").Value = -49691.5
$ws.Range("H81").Value = 3375.25
$ws.Range("I81").Value = 3375.25
$ws.Range("K81").Value = 6750.5
$ws.Range("M81").Value = -5689.5
$ws.Range("H84").Value = 3375.25
$ws.Range("I84").Value = 3375.25
$ws.Range("K84").Value = 33752.5
$ws.Range("M84").Value = -28448.5
$ws.Range("H122").Value = 1346
$ws.Range("I122").Value = 1087.5
$ws.Range("K122").Value = 3262.5
$ws.Range("M122").Value = -812.5
$ws.Range("H126").Value = 2958.2
$ws.Range("I126").Value = 1948
$ws.Range("K126").Value = 5844
$ws.Range("M126").Value = -3374
